# "agrego cv en word"
#
# The "Información Básica" block on the resume uses a negative left
# indent (-1134 twips / -56.7 pt) on its label/value paragraphs. This
# edit tightens that indent to -993 twips (-49.65 pt) for the ten
# paragraphs that make up that block (AGE, EMAIL, TELÉFONO, DIRECCIÓN
# and IDIOMAS labels plus their corresponding values), leaving the
# section heading ("Información Básica"), the trailing blank line and
# every other indent in the document untouched.

$d = $word.ActiveDocument

# Twips -> points: Word's ParagraphFormat.LeftIndent is expressed in
# points, and 1 point = 20 twips, so -993 twips == -49.65 pt.
$newIndentPt = -49.65

# Each value is unique within the document, so Find can locate the
# exact paragraph that owns it regardless of its position/index.
$targets = @(
    "AGE:",
    "39",
    "EMAIL:",
    "mauriciodiaz.xx@gmail.com",
    "TELÉFONO:",
    "+5491153117489",
    "DIRECCIÓN:",
    "San Pedro 170, Sarandi (Avellaneda) Buenos Aires, Argentina",
    "IDIOMAS:",
    "Inglés, Español"
)

foreach ($t in $targets) {
    $range = $d.Content
    $found = $range.Find.Execute($t, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
    if ($found) {
        $para = $range.Paragraphs.Item(1)
        $para.Format.LeftIndent = $newIndentPt
    }
}
